$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns
$ws.Range("H1").Value = "Tempo Heuristica"
$ws.Range("I1").Value = "Tempo Total"

# Update Row 2
$ws.Range("G2").Value = 0.02124500274658203
$ws.Range("H2").Value = 0.003348112106323242
$ws.Range("I2").Value = 0.02459311485290527

# Update Row 3
$ws.Range("G3").Value = 0.01173210144042969
$ws.Range("H3").Value = 0.002631664276123047
$ws.Range("I3").Value = 0.01436376571655273

# Update Row 4
$ws.Range("G4").Value = 0.01089954376220703
$ws.Range("H4").Value = 0.001598358154296875
$ws.Range("I4").Value = 0.01249790191650391
